$wb = $excel.ActiveWorkbook

$mdDisplay = "3c20df38-ca64-4e6d-a0de-282910fbaba4.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef28b28d3e96bb23d89bc076426ce4a525a65553/e2e/3c20df38-ca64-4e6d-a0de-282910fbaba4.md"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a25c45614ed5bcf67a468fca0a5626188a23ad3/e2e/3c20df38-ca64-4e6d-a0de-282910fbaba4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef28b28d3e96bb23d89bc076426ce4a525a65553/e2e/3c20df38-ca64-4e6d-a0de-282910fbaba4.md."

# ---- zh-cn sheet, row 7 ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "3c20df38-ca64-4e6d-a0de-282910fbaba4.f0dcec2c18ca77bcfe6a52a5e32c897d991a7197.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-20 17:03:25"
$wsZh.Range("P7").Value = $errorMsg

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276

# ---- de-de sheet, row 7 ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "3c20df38-ca64-4e6d-a0de-282910fbaba4.f0dcec2c18ca77bcfe6a52a5e32c897d991a7197.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-20 17:03:32"
$wsDe.Range("P7").Value = $errorMsg

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay)
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276
